$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.231.59'
$ws.Range("E2").Value = '  -0.61%  '

# Row 3
$ws.Range("D3").Value = '1.831.40'
$ws.Range("E3").Value = '  -0.60%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.37%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.03%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6003'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.42%  '

# Row 7
$ws.Range("E7").Value = '  +0.29%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.06993'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.49%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2760'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.60%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.31'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.45%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07616'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.38%  '

# Row 12
$ws.Range("D12").Value = '1.837.56'
$ws.Range("E12").Value = '  -0.34%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.760'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.12%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6288'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.13%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000009730'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.24%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '78.52'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.92%  '

# Row 17
$ws.Range("D17").Value = '28.680.03'
$ws.Range("E17").Value = '  -2.64%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.693'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -8.88%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '220.92'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.53%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.004'
$ws.Range("D20").Style = "Normal"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.56'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.94%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.871'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.84%  '

# Row 23
$ws.Range("E23").Value = '  +0.42%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '156.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.50%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.966'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.00%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1289'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.09%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.54'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.44%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.455'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.16%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.06429'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -11.34%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.439'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.86%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.841'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.69%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.754'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.93%  '

# Row 33
$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.727'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.98%  '

# Row 34
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.092'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.78%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6459'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.57%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.538'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.73%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.734'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.13%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01751'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.65%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.586'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.04%  '

# Row 40
$ws.Range("D40").Value = '1.170.94'
$ws.Range("E40").Value = '  -5.03%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8906'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.71%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.004'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.27%  '

# Row 43
$ws.Range("D43").Value = '1.983.88'
$ws.Range("E43").Value = '  -1.05%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.26'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.90%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.16'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.83%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000112'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.77%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05585'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.24%  '

# Row 48
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.589'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.43%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.464'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.68%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4556'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.39%  '

# Row 51
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.385'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -8.15%  '
